$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.868.72"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "2.715.82"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("D5").Value = "600.58"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").Value = "162.72"
$ws.Range("E6").Value = "  +3.85%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").Value = "2.715.18"
$ws.Range("E9").Value = "  +2.63%  "

$ws.Range("E10").Value = "  +0.52%  "

$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("E13").Value = "  +3.31%  "

$ws.Range("D14").Value = "28.53"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("D15").Value = "3.194.14"
$ws.Range("E15").Value = "  +2.12%  "

$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "68.739.06"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").Value = "2.712.78"
$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("E19").Value = "  +4.48%  "

$ws.Range("D20").Value = "7.69"
$ws.Range("E20").Value = "  +4.81%  "

$ws.Range("D21").Value = "366.17"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("E22").Value = "  +3.01%  "

$ws.Range("E23").Value = "  +2.89%  "

$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +2.96%  "

$ws.Range("D25").Value = "74.23"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "9.93"
$ws.Range("E27").Value = "  +1.72%  "

$ws.Range("D28").Value = "2.839.31"
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").Value = "594.59"
$ws.Range("E30").Value = "  +6.39%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").Value = "8.32"
$ws.Range("E32").Value = "  +3.41%  "

$ws.Range("D33").Value = "1.46"
$ws.Range("E33").Value = "  +3.47%  "

$ws.Range("E34").Value = "  +5.05%  "

$ws.Range("E35").Value = "  +3.55%  "

$ws.Range("D36").Value = "1.64"
$ws.Range("E36").Value = "  +5.72%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "19.93"
$ws.Range("E38").Value = "  +1.31%  "

$ws.Range("D39").Value = "160.29"
$ws.Range("E39").Value = "  -0.84%  "

$ws.Range("D40").Value = "0.381"
$ws.Range("E40").Value = "  +2.63%  "

$ws.Range("E41").Value = "  +2.67%  "

$ws.Range("D42").Value = "5.46"
$ws.Range("E42").Value = "  +2.72%  "

$ws.Range("D43").Value = "2.71"
$ws.Range("E43").Value = "  +4.24%  "

$ws.Range("E44").Value = "  +1.20%  "

$ws.Range("D45").Value = "0.0₆0319"
$ws.Range("E45").Value = "  -4.80%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").Value = "158.75"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  +6.17%  "

$ws.Range("E49").Value = "  +5.83%  "

$ws.Range("D50").Value = "0.607"
$ws.Range("E50").Value = "  +7.55%  "

$ws.Range("D51").Value = "22.21"
$ws.Range("E51").Value = "  +0.71%  "
